$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

# Update the "SÜRE" (duration) values in column B from 23 to 27 for DateNumber_1/2
$ws.Range("B8").Value = 27
$ws.Range("B9").Value = 27
$ws.Range("B11").Value = 27
$ws.Range("B12").Value = 27

# Fill in previously empty duration values for DateNumber_3/4
$ws.Range("B14").Value = 27
$ws.Range("B15").Value = 27
$ws.Range("B17").Value = 27
$ws.Range("B18").Value = 27

# Move the active selection to B16 as in the authored workbook
$ws.Activate()
$ws.Range("B16").Select()
